$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rework the sentence about the Product Backlog / user-stories paragraph.
#    Original: "final product. They can also be termed as user-stories. "
#    New:      "final product. They are derived from user-stories. The list
#               is prioritised by importance and story points denote the
#               weight and complexity of the task. "
#    This also relocates the hidden "_GoBack" bookmark into the middle of
#    the new sentence (right after "is prioritised by ").
# ---------------------------------------------------------------------------

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Text = "final product. They can also be termed as user-stories. "
$found = $r.Find.Execute()

if ($found) {
    $full = "final product. They are derived from user-stories. The list is prioritised by importance and story points denote the weight and complexity of the task. "
    $r.Text = $full
    $baseStart = $r.Start

    # Character offsets (relative to baseStart) of each logical chunk that
    # needs to become its own run in the final XML.
    $bounds = @(
        @(0, 20),    # "final product. They "
        @(20, 36),   # "are derived from"
        @(36, 51),   # " user-stories. "
        @(51, 60),   # "The list "
        @(60, 78),   # "is prioritised by "
        @(78, 152)   # "importance and story points denote the weight and complexity of the task. "
    )

    foreach ($b in $bounds) {
        $segStart = $baseStart + $b[0]
        $segEnd = $baseStart + $b[1]
        $seg = $d.Range($segStart, $segEnd)
        # Toggle bold on/off to force this chunk of text onto its own run
        # without actually changing any visible formatting.
        $seg.Font.Bold = 1
        $seg.Font.Bold = 0
    }

    # Move the "_GoBack" bookmark so that it now sits right after
    # "is prioritised by " (i.e. between that run and the following
    # "importance and story points..." run). Adding a bookmark with a name
    # that already exists elsewhere in the document relocates it.
    $bookmarkPos = $baseStart + 78
    $bmRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
